# Update evaluation metrics across the three worksheets with the final
# evaluation results.

$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.297153024911032
$wsSummary.Range("C2").Value = 0.06619385342789598
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.1241685144124169
$wsSummary.Range("F2").Value = 0.2616822429906542
$wsSummary.Range("G2").Value = 0.6482635796972396
$wsSummary.Range("H2").Value = 0.7215088282504013
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 395
$wsSummary.Range("K2").Value = 139
$wsSummary.Range("L2").Value = 0

# --- Sheet "Classification Report" ---
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$wsReport.Range("B2").Value = 1
$wsReport.Range("C2").Value = 0.2602996254681648
$wsReport.Range("D2").Value = 0.413075780089153

# Row 3 ("1")
$wsReport.Range("B3").Value = 0.06619385342789598
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.1241685144124169

# Row 4 ("accuracy")
$wsReport.Range("B4").Value = 0.297153024911032
$wsReport.Range("C4").Value = 0.297153024911032
$wsReport.Range("D4").Value = 0.297153024911032
$wsReport.Range("E4").Value = 0.297153024911032

# Row 5 ("macro avg")
$wsReport.Range("B5").Value = 0.533096926713948
$wsReport.Range("C5").Value = 0.6301498127340824
$wsReport.Range("D5").Value = 0.2686221472507849

# Row 6 ("weighted avg")
$wsReport.Range("B6").Value = 0.9534758503487208
$wsReport.Range("C6").Value = 0.297153024911032
$wsReport.Range("D6").Value = 0.3986818237920914

# --- Sheet "Confusion Matrix" ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$wsConfusion.Range("B2").Value = 139
$wsConfusion.Range("C2").Value = 395

# Row 3 ("Actual 1")
$wsConfusion.Range("B3").Value = 0
$wsConfusion.Range("C3").Value = 28
